$wb = $excel.ActiveWorkbook

# Add a new worksheet named "18" after the last existing sheet
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "18"

# Match the page-margin layout used by the other data sheets
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# Populate the random-mem-grav sample rows (formula label in column A, values B:Q)
$ws.Range("A1").Value = "memProtectedDiv(cos(cos(0)), sub(sin(sin(sin(sin(read(a0, a2))))), write(a0, abs(0), a2)))"
$ws.Range("B1").Value = -53.31
$ws.Range("C1").Value = -57.58
$ws.Range("D1").Value = -102.85
$ws.Range("E1").Value = -583.55
$ws.Range("F1").Value = -590.15
$ws.Range("G1").Value = -540
$ws.Range("H1").Value = -460
$ws.Range("I1").Value = -372.36
$ws.Range("J1").Value = -297.4
$ws.Range("K1").Value = -347.22
$ws.Range("L1").Value = -321.09
$ws.Range("M1").Value = -664.26
$ws.Range("N1").Value = -766.38
$ws.Range("O1").Value = -1101.66
$ws.Range("P1").Value = -1474.13
$ws.Range("Q1").Value = -1690.89

$ws.Range("A2").Value = "memProtectedDiv(memProtectedDiv(a2, conditional(sin(a1), 0)), protectedLog(cos(write(a0, a2, sin(sin(add(read(a0, 0), a2)))))))"
$ws.Range("B2").Value = -298.07
$ws.Range("C2").Value = -268.19
$ws.Range("D2").Value = -179.47
$ws.Range("E2").Value = -200.48
$ws.Range("F2").Value = -856.99
$ws.Range("G2").Value = -796.25
$ws.Range("H2").Value = -816.22
$ws.Range("I2").Value = -918.39
$ws.Range("J2").Value = -1041.32
$ws.Range("K2").Value = -1206.02
$ws.Range("L2").Value = -1225.25
$ws.Range("M2").Value = -1296.17
$ws.Range("N2").Value = -1394
$ws.Range("O2").Value = -1396.53
$ws.Range("P2").Value = -1378.41
$ws.Range("Q2").Value = -1731.88

$ws.Range("A3").Value = "memProtectedDiv(add(add(a2, a2), read(a0, a1)), conditional(sin(sin(sin(sin(read(a0, cos(protectedLog(memProtectedDiv(a1, a2)))))))), write(a0, limit(sub(a1, a1), protectedLog(conditional(sin(sin(read(a0, 0))), write(a0, sin(add(a2, a2)), sub(protectedLog(read(a0, 0)), write(a0, 0, cos(conditional(a1, a2))))))), a1), a2)))"
$ws.Range("B3").Value = -88.03
$ws.Range("C3").Value = -76.05
$ws.Range("D3").Value = -68.04000000000001
$ws.Range("E3").Value = -129.26
$ws.Range("F3").Value = -651.55
$ws.Range("G3").Value = -547.71
$ws.Range("H3").Value = -427.1
$ws.Range("I3").Value = -376.66
$ws.Range("J3").Value = -414.86
$ws.Range("K3").Value = -387.87
$ws.Range("L3").Value = -434.88
$ws.Range("M3").Value = -796.75
$ws.Range("N3").Value = -969.4
$ws.Range("O3").Value = -1125.85
$ws.Range("P3").Value = -1421.98
$ws.Range("Q3").Value = -1678.01

$ws.Range("A4").Value = "sub(sub(sub(read(a0, 0), memProtectedDiv(a2, abs(read(a0, sub(cos(protectedLog(add(0, a1))), a1))))), write(a0, sin(read(a0, memProtectedDiv(a2, a1))), a2)), a2)"
$ws.Range("B4").Value = -157.45
$ws.Range("C4").Value = -122.33
$ws.Range("D4").Value = -110.08
$ws.Range("E4").Value = -182.15
$ws.Range("F4").Value = -675.09
$ws.Range("G4").Value = -771.4400000000001
$ws.Range("H4").Value = -846.36
$ws.Range("I4").Value = -894.27
$ws.Range("J4").Value = -1034.44
$ws.Range("K4").Value = -1097.06
$ws.Range("L4").Value = -1186.2
$ws.Range("M4").Value = -1255.08
$ws.Range("N4").Value = -1269.79
$ws.Range("O4").Value = -1316.21
$ws.Range("P4").Value = -1345.97
$ws.Range("Q4").Value = -1634.54

$ws.Range("A5").Value = "memProtectedDiv(sub(write(a0, a1, sub(sub(sub(read(a0, 0), a2), a2), abs(protectedLog(0)))), write(a0, conditional(cos(a1), a1), add(write(a0, a1, add(write(a0, a1, add(a2, add(0, a2))), add(limit(a2, 0, write(a0, a1, 0)), a2))), a2))), abs(a2))"
$ws.Range("B5").Value = -63.25
$ws.Range("C5").Value = -63.72
$ws.Range("D5").Value = -58.33
$ws.Range("E5").Value = -106.62
$ws.Range("F5").Value = -371.9
$ws.Range("G5").Value = -386.74
$ws.Range("H5").Value = -493.96
$ws.Range("I5").Value = -638.3200000000001
$ws.Range("J5").Value = -870.0700000000001
$ws.Range("K5").Value = -844.96
$ws.Range("L5").Value = -982.09
$ws.Range("M5").Value = -957.25
$ws.Range("N5").Value = -928.6799999999999
$ws.Range("O5").Value = -355.68
$ws.Range("P5").Value = -616.05
$ws.Range("Q5").Value = -1483.02

$ws.Range("A6").Value = "sub(sub(0, memProtectedDiv(conditional(a2, read(a0, a2)), write(a0, 0, sin(a2)))), memProtectedDiv(conditional(read(a0, 0), 0), a2))"
$ws.Range("B6").Value = -68.2
$ws.Range("C6").Value = -77.67
$ws.Range("D6").Value = -68.23
$ws.Range("E6").Value = -101.44
$ws.Range("F6").Value = -422.63
$ws.Range("G6").Value = -382.93
$ws.Range("H6").Value = -373.76
$ws.Range("I6").Value = -365.42
$ws.Range("J6").Value = -487.3
$ws.Range("K6").Value = -624.72
$ws.Range("L6").Value = -714.51
$ws.Range("M6").Value = -789.6799999999999
$ws.Range("N6").Value = -771.95
$ws.Range("O6").Value = -594
$ws.Range("P6").Value = -1362.07
$ws.Range("Q6").Value = -1493.48

$ws.Range("A7").Value = "sub(sub(sub(sub(read(a0, add(add(memProtectedDiv(0, 0), a1), a1)), write(a0, limit(0, sub(add(0, a2), protectedLog(a1)), abs(a2)), add(add(memProtectedDiv(protectedLog(abs(a2)), sin(a2)), memProtectedDiv(a1, conditional(a2, a1))), a2))), sin(read(a0, add(read(a0, a2), a2)))), a2), add(a2, a2))"
$ws.Range("B7").Value = -87.09
$ws.Range("C7").Value = -88.98999999999999
$ws.Range("D7").Value = -84.65000000000001
$ws.Range("E7").Value = -189.07
$ws.Range("F7").Value = -210.36
$ws.Range("G7").Value = -219.76
$ws.Range("H7").Value = -246.34
$ws.Range("I7").Value = -244.69
$ws.Range("J7").Value = -338.3
$ws.Range("K7").Value = -422.06
$ws.Range("L7").Value = -457.37
$ws.Range("M7").Value = -538.99
$ws.Range("N7").Value = -921.0700000000001
$ws.Range("O7").Value = -1432.88
$ws.Range("P7").Value = -1572.6
$ws.Range("Q7").Value = -1684.28

$ws.Range("A8").Value = "sub(read(a0, a2), memProtectedDiv(abs(memProtectedDiv(memProtectedDiv(write(a0, protectedLog(0), add(a2, a2)), a1), a1)), a2))"
$ws.Range("B8").Value = -254.59
$ws.Range("C8").Value = -305.07
$ws.Range("D8").Value = -251.89
$ws.Range("E8").Value = -306.78
$ws.Range("F8").Value = -803.49
$ws.Range("G8").Value = -902.97
$ws.Range("H8").Value = -1105.06
$ws.Range("I8").Value = -1204.28
$ws.Range("J8").Value = -1313.37
$ws.Range("K8").Value = -1460.39
$ws.Range("L8").Value = -1477.73
$ws.Range("M8").Value = -1497.64
$ws.Range("N8").Value = -1543.71
$ws.Range("O8").Value = -1529.84
$ws.Range("P8").Value = -1460.22
$ws.Range("Q8").Value = -1714.57

$ws.Range("A9").Value = "sub(read(a0, a2), add(a2, add(add(a2, write(a0, sin(a1), sub(a1, memProtectedDiv(conditional(memProtectedDiv(a1, a2), a2), a2)))), add(a2, cos(read(a0, 0))))))"
$ws.Range("B9").Value = -99.48999999999999
$ws.Range("C9").Value = -92.16
$ws.Range("D9").Value = -97.34
$ws.Range("E9").Value = -180.29
$ws.Range("F9").Value = -194.12
$ws.Range("G9").Value = -205.33
$ws.Range("H9").Value = -251.55
$ws.Range("I9").Value = -278.61
$ws.Range("J9").Value = -351.4
$ws.Range("K9").Value = -402.72
$ws.Range("L9").Value = -483.8
$ws.Range("M9").Value = -627.42
$ws.Range("N9").Value = -1101.22
$ws.Range("O9").Value = -1375.39
$ws.Range("P9").Value = -1565.57
$ws.Range("Q9").Value = -1757.59

# Select A1 on the new sheet and make it the active tab, matching the target workbook view
$ws.Range("A1").Select()
$ws.Activate()
